$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 191
$ws.Range("F3").Value = 1025
$ws.Range("F7").Value = 724
$ws.Range("F8").Value = 13071
$ws.Range("F9").Value = 2265
$ws.Range("B10").Value = "'2024-10-02"
$ws.Range("C10").Value = "建德·逆光ZERO动漫游戏展（取消）"
$ws.Range("D10").Value = "南山路1号 杭州新安雷迪森酒店"
$ws.Range("E10").Value = "2024.10.02 10:00-10.02 17:00"
$ws.Range("F10").Value = 7
$ws.Range("G10").Value = "不可售"
$ws.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=91698"
$ws.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202408/vQITgOEI1724885714305.jpeg"
$ws.Range("F11").Value = 297
$ws.Range("F12").Value = 53668
$ws.Range("G12").Value = "暂时售罄"
$ws.Range("F13").Value = 1291
$ws.Range("F14").Value = 298
$ws.Range("F15").Value = 299
$ws.Range("G15").Value = 128
$ws.Range("F20").Value = 844
$ws.Range("F21").Value = 5076
$ws.Range("F22").Value = 1238
$ws.Range("F25").Value = 33
$ws.Range("F27").Value = 47
$ws.Range("F28").Value = 1177
$ws.Range("F29").Value = 77
$ws.Range("F30").Value = 23
$ws.Range("F31").Value = 140
$ws.Range("F32").Value = 319
$ws.Range("F35").Value = 56
$ws.Range("F36").Value = 37
$ws.Range("F37").Value = 4664
$ws.Range("F39").Value = 4711
$ws.Range("F40").Value = 8669
$ws.Range("F43").Value = 110
$ws.Range("F44").Value = 198
$ws.Range("F45").Value = 399
$ws.Range("F46").Value = 95
$ws.Range("F47").Value = 65
$ws.Range("F48").Value = 4153
$ws.Range("F49").Value = 171

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 122
$ws.Range("F12").Value = 1104
$ws.Range("F13").Value = 7
$ws.Range("F20").Value = 84

# ---- Sheet: 本地生活 ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 531
$ws.Range("F4").Value = 136
$ws.Range("G4").Value = "已售罄"
$ws.Range("F5").Value = 28

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 531
$ws.Range("F4").Value = 191
$ws.Range("F5").Value = 1026
$ws.Range("F7").Value = 724
$ws.Range("F8").Value = 13072
$ws.Range("F9").Value = 13072
$ws.Range("F10").Value = 2265
$ws.Range("B11").Value = "'2024-10-03"
$ws.Range("C11").Value = "杭州·COMICUP30"
$ws.Range("D11").Value = "阳城路雅澳杭州电商产业园西侧约200米 杭州大会展中心"
$ws.Range("E11").Value = "2024.10.03 10:00-10.07 16:00"
$ws.Range("F11").Value = 53668
$ws.Range("G11").Value = "暂时售罄"
$ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=92785"
$ws.Range("I11").Value = "//i0.hdslb.com/bfs/openplatform/202409/UIz8Mjqi1727086527159.png"
$ws.Range("F12").Value = 298
$ws.Range("F17").Value = 844
$ws.Range("F19").Value = 5076
$ws.Range("F20").Value = 1238
$ws.Range("F21").Value = 28
$ws.Range("F22").Value = 122
$ws.Range("F24").Value = 48
$ws.Range("F26").Value = 1177
$ws.Range("F28").Value = 77
$ws.Range("F29").Value = 23
$ws.Range("F30").Value = 140
$ws.Range("F32").Value = 319
$ws.Range("F34").Value = 37
$ws.Range("F35").Value = 4664
$ws.Range("F36").Value = 4711
$ws.Range("F37").Value = 8669
$ws.Range("F40").Value = 198
$ws.Range("F43").Value = 95
$ws.Range("F44").Value = 4153
$ws.Range("F46").Value = 84
